$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = "Dempsey Roll"
$ws.Range("E19").Value = "something that probably lowers defense more and more each round, but leads to higher damage with each round? And allows for an attack each round? Risk reward thing"

# Row 20
$ws.Range("A20").Value = "The Final Nail"
$ws.Range("B20").Value = "Enemy"
$ws.Range("C20").Value = "Instant"
$ws.Range("E20").Value = "use a mark to do something cool"
$ws.Range("G20").Value = "And the hammer that drove it"
$ws.Range("H20").Value = "??"
$ws.Range("I20").Value = "??"
$ws.Range("K20").Value = "something also with marks?"

# Row 21
$ws.Range("A21").Value = "Tight 5"
$ws.Range("E21").Value = "Something using charisma"

# Row 22
$ws.Range("A22").Value = "And the winner is…"

# Row 23
$ws.Range("E23").Value = "enchant a small container of liquid. When broken over a weapon, the next attack made with that weapon does an additional 1d8 fire damage"
$ws.Range("F23").Value = "exhaust"
$ws.Range("A23").Value = "Flame Oil"
$ws.Range("B23").Value = "Self"
$ws.Range("C23").Value = "X Rnds"
$ws.Range("D23").Value = "X = Level"

# Update selection to match the active cell state captured in the diff
$ws.Range("E24").Select()
